# Insert a new weekly record at row 80 for "Agrícola del Norte S.A. de Arica -
# Cebollín baby", pushing the previously-recorded rows (old 80..114) down to
# (81..115) and giving the sheet a brand-new first row with the latest
# price-report figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 80; Excel shifts rows 80-114 (and
# their formatting) down to 81-115 and grows the used range to A1:R115
# automatically.
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new observation.
$ws.Range("A80").Value = 1
$ws.Range("B80").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C80").Value = "Arica y Parinacota"
$ws.Range("D80").Value = 44917
$ws.Range("E80").Value = 15
$ws.Range("F80").Value = 100112038
$ws.Range("G80").Value = "Cebollín baby"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 250
$ws.Range("K80").Value = 5500
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = 5800
$ws.Range("N80").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O80").Value = "Región de Arica y Parinacota"
$ws.Range("P80").Value = 2900
$ws.Range("Q80").Value = 2
$ws.Range("R80").Value = "Hortaliza"
